$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: standalone note
$ws.Range("A10").Value = "*DON'T NEED TO WORRY ABOUT INCORRECT VS CORRECT"

# Contrast names (rows 12-16), entered in the same order the original
# author used so the shared-string table's insertion order matches.
$ws.Range("A13").Value = "2_Risk>Neutral"
$ws.Range("A15").Value = "4_NoGoRisk>GoRisk"
$ws.Range("A14").Value = "3_RiskNoGo.v.Go>NeutralNoGo.v.Go"
$ws.Range("A16").Value = "5_NoGoNeutral>GoNeutral"

# Averaging note, repeated across B12:B16
$ws.Range("B12").Value = "average pre, average post, average across, pre-post change, "

# Row 18: closing note
$ws.Range("A18").Value = "*At group level, you can flip contrasts (e.g. look at Go>NoGo as sanity check for motor activity)"

# A12 reuses the existing "1_NoGo>Go" shared string
$ws.Range("A12").Value = "1_NoGo>Go"

$ws.Range("B13").Value = "average pre, average post, average across, pre-post change, "
$ws.Range("B14").Value = "average pre, average post, average across, pre-post change, "
$ws.Range("B15").Value = "average pre, average post, average across, pre-post change, "
$ws.Range("B16").Value = "average pre, average post, average across, pre-post change, "

# Column A width adjustment (per diff: 23.1640625 -> 32.6640625, bestFit removed)
$ws.Columns("A").ColumnWidth = 31.83

# View settings: zoom 150%, selection at A20
$excel.ActiveWindow.Zoom = 150
$ws.Range("A20").Select()
